$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 73.111115
$ws.Range("I55").Value = 73.111115
$ws.Range("K55").Value = 73.111115
$ws.Range("M55").Value = 140.888885

$ws.Range("H129").Value = 1010.90564
$ws.Range("I129").Value = 448.5
$ws.Range("J129").Value = 1032.9608
$ws.Range("K129").Value = 1345.5
$ws.Range("L129").Value = 3098.8824
$ws.Range("M129").Value = 3654.5
$ws.Range("N129").Value = -13098.8824

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1309.2667
$ws.Range("I2").Value = 1406.7
$ws.Range("J2").Value = 1114.4
$ws.Range("K2").Value = 1406.7
$ws.Range("L2").Value = 1114.4
$ws.Range("M2").Value = -1293.7
$ws.Range("N2").Value = -1340.4

$ws.Range("H32").Value = 6770.55
$ws.Range("I32").Value = 5386.0547
$ws.Range("K32").Value = 5386.0547
$ws.Range("M32").Value = -5099.0547

$ws.Range("H63").Value = 2005
$ws.Range("I63").Value = 2005
$ws.Range("K63").Value = 2005
$ws.Range("M63").Value = -1319

$ws.Range("H66").Value = 2005
$ws.Range("I66").Value = 2005
$ws.Range("K66").Value = 10025
$ws.Range("M66").Value = -6593

$ws.Range("H92").Value = 23800
$ws.Range("J92").Value = 23800
$ws.Range("L92").Value = 23800
$ws.Range("N92").Value = -28792

$ws.Range("H102").Value = 2471105.8
$ws.Range("I102").Value = 2647313.2
$ws.Range("J102").Value = 4200
$ws.Range("K102").Value = 2647313.2
$ws.Range("L102").Value = 4200
$ws.Range("M102").Value = -2645691.2
$ws.Range("N102").Value = -7444

$ws.Range("H116").Value = 1309.2667
$ws.Range("I116").Value = 1406.7
$ws.Range("J116").Value = 1114.4
$ws.Range("K116").Value = 1406.7
$ws.Range("L116").Value = 1114.4
$ws.Range("M116").Value = 887.3
$ws.Range("N116").Value = -5702.4

$ws.Range("H128").Value = 67711.60000000001
$ws.Range("J128").Value = 67711.60000000001
$ws.Range("L128").Value = 67711.60000000001
$ws.Range("N128").Value = -77671.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1309.2667
$ws.Range("I3").Value = 1406.7
$ws.Range("J3").Value = 1114.4
$ws.Range("K3").Value = 1406.7
$ws.Range("L3").Value = 1114.4
$ws.Range("M3").Value = -1292.7
$ws.Range("N3").Value = -1342.4

$ws.Range("H80").Value = 223.2
$ws.Range("J80").Value = 237.38095
$ws.Range("L80").Value = 237.38095
$ws.Range("N80").Value = -2233.38095

$ws.Range("H83").Value = 223.2
$ws.Range("J83").Value = 237.38095
$ws.Range("L83").Value = 1186.90475
$ws.Range("N83").Value = -11170.90475

$ws.Range("H94").Value = 1093.8125
$ws.Range("I94").Value = 992.5599999999999
$ws.Range("J94").Value = 1455.4286
$ws.Range("K94").Value = 992.5599999999999
$ws.Range("L94").Value = 1455.4286
$ws.Range("M94").Value = -541.5599999999999
$ws.Range("N94").Value = -2357.4286

$ws.Range("H104").Value = 48000
$ws.Range("J104").Value = 48000
$ws.Range("L104").Value = 48000
$ws.Range("N104").Value = -54988

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 8750
$ws.Range("I32").Value = 8750
$ws.Range("K32").Value = 8750
$ws.Range("M32").Value = -8434

$ws.Range("H103").Value = 6233
$ws.Range("I103").Value = 6233
$ws.Range("K103").Value = 6233
$ws.Range("M103").Value = -5061

$ws.Range("H132").Value = 3281.3333
$ws.Range("I132").Value = 2636.5
$ws.Range("J132").Value = 4571
$ws.Range("K132").Value = 7909.5
$ws.Range("L132").Value = 13713
$ws.Range("M132").Value = -5379.5
$ws.Range("N132").Value = -18773

$ws.Range("H134").Value = 3965.775
$ws.Range("I134").Value = 3084.0527
$ws.Range("K134").Value = 9252.158100000001
$ws.Range("M134").Value = -6717.158100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 33430234
$ws.Range("I9").Value = 2000
$ws.Range("J9").Value = 41787292
$ws.Range("K9").Value = 6000
$ws.Range("L9").Value = 125361876
$ws.Range("M9").Value = -5776
$ws.Range("N9").Value = -125362324

$ws.Range("H14").Value = 1361.64
$ws.Range("I14").Value = 1361.64
$ws.Range("K14").Value = 4084.92
$ws.Range("M14").Value = -3911.92

$ws.Range("H34").Value = 2335.1562
$ws.Range("I34").Value = 225
$ws.Range("J34").Value = 3160.8696
$ws.Range("K34").Value = 675
$ws.Range("L34").Value = 9482.6088
$ws.Range("M34").Value = -591
$ws.Range("N34").Value = -9650.6088

$ws.Range("H69").Value = 1467.2413
$ws.Range("J69").Value = 1622
$ws.Range("L69").Value = 4866
$ws.Range("N69").Value = -6488

$ws.Range("H72").Value = 1467.2413
$ws.Range("J72").Value = 1622
$ws.Range("L72").Value = 14598
$ws.Range("N72").Value = -22710

$ws.Range("H80").Value = 2923.0667
$ws.Range("I80").Value = 4450
$ws.Range("J80").Value = 2688.1538
$ws.Range("K80").Value = 13350
$ws.Range("L80").Value = 8064.4614
$ws.Range("M80").Value = -12414
$ws.Range("N80").Value = -9936.4614

$ws.Range("H83").Value = 2923.0667
$ws.Range("I83").Value = 4450
$ws.Range("J83").Value = 2688.1538
$ws.Range("K83").Value = 40050
$ws.Range("L83").Value = 24193.3842
$ws.Range("M83").Value = -35370
$ws.Range("N83").Value = -33553.3842

$ws.Range("H98").Value = 385.13635
$ws.Range("I98").Value = 396.2
$ws.Range("J98").Value = 361.42856
$ws.Range("K98").Value = 1188.6
$ws.Range("L98").Value = 1084.28568
$ws.Range("M98").Value = 309.4000000000001
$ws.Range("N98").Value = -4080.28568

$ws.Range("H129").Value = 2188.3333
$ws.Range("J129").Value = 2011.4286
$ws.Range("L129").Value = 6034.2858
$ws.Range("N129").Value = -16034.2858

$ws.Range("H130").Value = 3795.5557
$ws.Range("I130").Value = 440
$ws.Range("J130").Value = 7990
$ws.Range("K130").Value = 1320
$ws.Range("L130").Value = 23970
$ws.Range("M130").Value = 3700
$ws.Range("N130").Value = -34010

$ws.Range("H131").Value = 2185.5
$ws.Range("I131").Value = 3993.75
$ws.Range("J131").Value = 980
$ws.Range("K131").Value = 11981.25
$ws.Range("L131").Value = 2940
$ws.Range("M131").Value = -6941.25
$ws.Range("N131").Value = -13020

$ws.Range("H140").Value = 2044.0625
$ws.Range("I140").Value = 1630.8695
$ws.Range("J140").Value = 3100
$ws.Range("K140").Value = 4892.6085
$ws.Range("L140").Value = 9300
$ws.Range("M140").Value = 287.3914999999997
$ws.Range("N140").Value = -19660

$ws.Range("H141").Value = 4041
$ws.Range("I141").Value = 3026.6667
$ws.Range("J141").Value = 4379.1113
$ws.Range("K141").Value = 9080.000100000001
$ws.Range("L141").Value = 13137.3339
$ws.Range("N141").Value = -23497.3339

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H135").Value = 57484.445
$ws.Range("J135").Value = 57484.445
$ws.Range("L135").Value = 57484.445
$ws.Range("N135").Value = -67624.44500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 73900
$ws.Range("J128").Value = 73900
$ws.Range("L128").Value = 73900
$ws.Range("N128").Value = -83860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 45000
$ws.Range("J40").Value = 45000
$ws.Range("L40").Value = 45000
$ws.Range("N40").Value = -45298

$ws.Range("H97").Value = 11572
$ws.Range("J97").Value = 11572
$ws.Range("L97").Value = 11572
$ws.Range("N97").Value = -13554

$ws.Range("H132").Value = 2068.739
$ws.Range("I132").Value = 1834.1
$ws.Range("J132").Value = 3633
$ws.Range("K132").Value = 5502.299999999999
$ws.Range("L132").Value = 10899
$ws.Range("M132").Value = -2972.299999999999
$ws.Range("N132").Value = -15959

$ws.Range("H136").Value = 4927.38
$ws.Range("I136").Value = 2825.9524
$ws.Range("J136").Value = 6449.1035
$ws.Range("K136").Value = 8477.8572
$ws.Range("L136").Value = 19347.3105
$ws.Range("M136").Value = -5927.8572
$ws.Range("N136").Value = -24447.3105
